$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.061.78'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '''  -2.07%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = '''1.669.66'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '''  -1.47%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = '''  -0.19%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = '''216.95'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '''  -1.19%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = '''0.5114'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '''  +0.45%  '
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = '''1.005'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '''  -0.18%  '
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = '''  +0.57%  '
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = '''0.06414'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '''  +2.08%  '
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = '''21.83'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '''  -1.36%  '
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = '''0.07442'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '''  +1.29%  '
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = '''1.689.34'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '''  -0.53%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = '''4.507'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '''  -0.14%  '
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = '''  +1.02%  '
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = '''0.000008594'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '''  +1.93%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = '''64.49'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '''  -1.40%  '
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = '''26.098.13'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '''  -2.01%  '
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = '''4.949'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '''  -0.66%  '
$ws.Range("E18").ClearFormats()
$ws.Range("E20").Value = '''  -1.67%  '
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = '''193.87'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '''  +4.18%  '
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = '''6.224'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '''  -0.27%  '
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = '''  -0.16%  '
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = '''144.94'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '''  +0.27%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = '''7.618'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '''  +1.74%  '
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = '''  +3.61%  '
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = '''15.75'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '''  -0.20%  '
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = '''0.06448'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '''  +14.02%  '
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = '''  -1.13%  '
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = '''  -1.19%  '
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = '''3.552'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '''  +1.41%  '
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = '''  +1.07%  '
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = '''1.649'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '''  -0.27%  '
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = '''  +0.17%  '
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = '''0.6114'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '''  +2.12%  '
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = '''2.369'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '''  +0.14%  '
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = '''  +0.31%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = '''6.260'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '''  +7.33%  '
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = '''0.01605'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '''  -0.60%  '
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = '''1.092.57'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '''  -0.98%  '
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = '''0.8633'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '''  +0.83%  '
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = '''  +0.57%  '
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = '''100.55'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '''  +1.12%  '
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = '''1.818.07'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '''  -1.88%  '
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = '''0.00000000114'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '''  +1.93%  '
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = '''  +0.03%  '
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = '''  +0.40%  '
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = '''8.085'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '''  -0.59%  '
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = '''0.05238'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '''  -0.02%  '
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = '''0.4286'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '''  -0.86%  '
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = '''6.041'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '''  +4.57%  '
$ws.Range("E51").ClearFormats()
